# "Eintritte und Kiosk für 1.4.24"
# Fill in the missing "Minimal Abzug" (col C) / "Abzug [%]" (col D) figures
# for the two screenings added on 1.4.24 (rows 15 and 17 of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verleiherabgaben")
$ws.Activate()

# Row 15 (1018.574 Suisa code, 1.4.24 entry) was missing the Abzug [%] value.
$ws.Range("D15").Value = 50

# Row 17 (same Suisa code, repeated on 1.4.24) was missing both the
# Minimal Abzug and Abzug [%] values.
$ws.Range("C17").Value = 150
$ws.Range("D17").Value = 50

# Move the selection to C16, matching where the author ended up.
$ws.Range("C16").Select()
